$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three more rows identical to row 4 (Kots / images/kotesh.png)
$ws.Range("A5").Value = "Kots"
$ws.Range("B5").Value = "images/kotesh.png"

$ws.Range("A6").Value = "Kots"
$ws.Range("B6").Value = "images/kotesh.png"

$ws.Range("A7").Value = "Kots"
$ws.Range("B7").Value = "images/kotesh.png"
